$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cell values (row-by-row financial data refresh) ---
$ws.Range("D2").Value = 45834
$ws.Range("E2").Value = 9411
$ws.Range("F2").Value = 9411
$ws.Range("G2").Value = 9191
$ws.Range("H2").Value = 8197
$ws.Range("I2").Value = 8098
$ws.Range("J2").Value = 99
$ws.Range("K2").Value = 840501
$ws.Range("L2").Value = 785861
$ws.Range("M2").Value = 54639
$ws.Range("N2").Value = 47657
$ws.Range("O2").Value = 6983
$ws.Range("P2").Value = 11719
$ws.Range("Q2").Value = -1305
$ws.Range("R2").Value = -12375
$ws.Range("S2").Value = 13664
$ws.Range("T2").Value = 1750
$ws.Range("V2").Value = 131016
$ws.Range("W2").Value = 20.53
$ws.Range("X2").Value = 17.88
$ws.Range("Y2").Value = 19.59
$ws.Range("Z2").Value = 1.25
$ws.Range("AA2").Value = 1438.27
$ws.Range("AB2").Value = 366.25
$ws.Range("AC2").Value = 3523
$ws.Range("AD2").Value = 3.96
$ws.Range("AE2").Value = 19153
$ws.Range("AF2").Value = 0.73
$ws.Range("AG2").Value = 188
$ws.Range("AH2").Value = 1.35
$ws.Range("AI2").Value = 5.79
$ws.Range("AJ2").Value = 248822827
$ws.Range("D3").Value = 51740
$ws.Range("E3").Value = 6962
$ws.Range("F3").Value = 6962
$ws.Range("G3").Value = 7001
$ws.Range("H3").Value = 5305
$ws.Range("I3").Value = 4855
$ws.Range("J3").Value = 450
$ws.Range("K3").Value = 902795
$ws.Range("L3").Value = 840889
$ws.Range("M3").Value = 61907
$ws.Range("N3").Value = 58862
$ws.Range("O3").Value = 3045
$ws.Range("P3").Value = 12797
$ws.Range("Q3").Value = 7464
$ws.Range("R3").Value = -8555
$ws.Range("S3").Value = 4361
$ws.Range("T3").Value = 1051
$ws.Range("V3").Value = 131317
$ws.Range("W3").Value = 13.46
$ws.Range("X3").Value = 10.25
$ws.Range("Y3").Value = 9.119999999999999
$ws.Range("Z3").Value = 0.61
$ws.Range("AA3").Value = 1358.31
$ws.Range("AB3").Value = 383.79
$ws.Range("AC3").Value = 1853
$ws.Range("AD3").Value = 4.54
$ws.Range("AE3").Value = 21665
$ws.Range("AF3").Value = 0.39
$ws.Range("AG3").Value = 141
$ws.Range("AH3").Value = 1.68
$ws.Range("AI3").Value = 7.91
$ws.Range("AJ3").Value = 271706454
$ws.Range("D4").Value = 49127
$ws.Range("E4").Value = 7123
$ws.Range("F4").Value = 7123
$ws.Range("G4").Value = 6884
$ws.Range("H4").Value = 5181
$ws.Range("I4").Value = 5016
$ws.Range("J4").Value = 165
$ws.Range("K4").Value = 934822
$ws.Range("L4").Value = 863950
$ws.Range("M4").Value = 70872
$ws.Range("N4").Value = 67826
$ws.Range("O4").Value = 3045
$ws.Range("P4").Value = 16297
$ws.Range("Q4").Value = -11258
$ws.Range("R4").Value = -3719
$ws.Range("S4").Value = 11978
$ws.Range("T4").Value = 1041
$ws.Range("V4").Value = 136111
$ws.Range("W4").Value = 14.5
$ws.Range("X4").Value = 10.55
$ws.Range("Y4").Value = 7.92
$ws.Range("Z4").Value = 0.5600000000000001
$ws.Range("AA4").Value = 1219.04
$ws.Range("AB4").Value = 334.89
$ws.Range("AC4").Value = 1568
$ws.Range("AD4").Value = 5.53
$ws.Range("AE4").Value = 20811
$ws.Range("AF4").Value = 0.42
$ws.Range("AG4").Value = 230
$ws.Range("AH4").Value = 2.65
$ws.Range("AI4").Value = 14.94
$ws.Range("AJ4").Value = 325935246
$ws.Range("D5").Value = 48098
$ws.Range("E5").Value = 5943
$ws.Range("F5").Value = 5943
$ws.Range("G5").Value = 5578
$ws.Range("H5").Value = 4250
$ws.Range("I5").Value = 4031
$ws.Range("J5").Value = 219
$ws.Range("K5").Value = 943500
$ws.Range("L5").Value = 866940
$ws.Range("M5").Value = 76560
$ws.Range("N5").Value = 70576
$ws.Range("O5").Value = 5985
$ws.Range("P5").Value = 16297
$ws.Range("Q5").Value = 15302
$ws.Range("R5").Value = -4273
$ws.Range("S5").Value = -9320
$ws.Range("T5").Value = 1513
$ws.Range("V5").Value = 127531
$ws.Range("W5").Value = 12.36
$ws.Range("X5").Value = 8.83
$ws.Range("Y5").Value = 5.83
$ws.Range("Z5").Value = 0.45
$ws.Range("AA5").Value = 1132.36
$ws.Range("AB5").Value = 369.8
$ws.Range("AC5").Value = 1237
$ws.Range("AD5").Value = 7.62
$ws.Range("AE5").Value = 21654
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 230
$ws.Range("AH5").Value = 2.44
$ws.Range("AI5").Value = 18.6
$ws.Range("AJ5").Value = 325935246
$ws.Range("D6").Value = 50554
$ws.Range("E6").Value = 7498
$ws.Range("F6").Value = 7498
$ws.Range("G6").Value = 7186
$ws.Range("H6").Value = 5381
$ws.Range("I6").Value = 5021
$ws.Range("K6").Value = 987936
$ws.Range("L6").Value = 903686
$ws.Range("M6").Value = 84250
$ws.Range("N6").Value = 76270
$ws.Range("P6").Value = 16297
$ws.Range("Q6").Value = 3544
$ws.Range("R6").Value = -8538
$ws.Range("S6").Value = 6364
$ws.Range("T6").Value = 1253
$ws.Range("V6").Value = 124314
$ws.Range("W6").Value = 14.83
$ws.Range("X6").Value = 10.64
$ws.Range("Y6").Value = 6.84
$ws.Range("Z6").Value = 0.5600000000000001
$ws.Range("AA6").Value = 1072.63
$ws.Range("AB6").Value = 416.99
$ws.Range("AC6").Value = 1540
$ws.Range("AD6").Value = 4.76
$ws.Range("AE6").Value = 23401
$ws.Range("AF6").Value = 0.31
$ws.Range("AI6").Value = 19.47
$ws.Range("AJ6").Value = 325935246
$ws.Range("E7").Value = 8074
$ws.Range("G7").Value = 8189
$ws.Range("H7").Value = 6136
$ws.Range("I7").Value = 5806
$ws.Range("K7").Value = 1045525
$ws.Range("L7").Value = 954356
$ws.Range("M7").Value = 91169
$ws.Range("N7").Value = 82244
$ws.Range("P7").Value = 17008
$ws.Range("Y7").Value = 7.32
$ws.Range("Z7").Value = 0.6
$ws.Range("AA7").Value = 1046.8
$ws.Range("AC7").Value = 1781
$ws.Range("AD7").Value = 3.79
$ws.Range("AE7").Value = 25234
$ws.Range("AF7").Value = 0.27
$ws.Range("AG7").Value = 349
$ws.Range("AH7").Value = 5.17
$ws.Range("AI7").Value = 19.6
$ws.Range("E8").Value = 7951
$ws.Range("G8").Value = 7857
$ws.Range("H8").Value = 5821
$ws.Range("I8").Value = 5378
$ws.Range("K8").Value = 1086348
$ws.Range("L8").Value = 990383
$ws.Range("M8").Value = 95965
$ws.Range("N8").Value = 86657
$ws.Range("P8").Value = 17008
$ws.Range("Y8").Value = 6.37
$ws.Range("Z8").Value = 0.55
$ws.Range("AA8").Value = 1032.03
$ws.Range("AC8").Value = 1650
$ws.Range("AD8").Value = 4.09
$ws.Range("AE8").Value = 26588
$ws.Range("AF8").Value = 0.25
$ws.Range("AG8").Value = 358
$ws.Range("AH8").Value = 5.31
$ws.Range("AI8").Value = 21.71
$ws.Range("E9").Value = 8405
$ws.Range("G9").Value = 8276
$ws.Range("H9").Value = 6093
$ws.Range("I9").Value = 5533
$ws.Range("K9").Value = 1122898
$ws.Range("L9").Value = 1022962
$ws.Range("M9").Value = 99935
$ws.Range("N9").Value = 91305
$ws.Range("P9").Value = 17185
$ws.Range("Y9").Value = 6.22
$ws.Range("Z9").Value = 0.55
$ws.Range("AA9").Value = 1023.63
$ws.Range("AC9").Value = 1698
$ws.Range("AD9").Value = 3.98
$ws.Range("AE9").Value = 28015
$ws.Range("AF9").Value = 0.24
$ws.Range("AG9").Value = 390
$ws.Range("AH9").Value = 5.78
$ws.Range("AI9").Value = 22.97

# --- Clear cells that no longer have data in the refreshed report ---
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()

